# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns for the zh-cn and de-de handback tables, updates the status text and widens the
# columns that now hold longer hyperlink text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/677d2ed9ce93812cf171508034449025ff985881/e2e/3b55afe9-fca1-4f9d-b841-af26fca2fc20.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/677d2ed9ce93812cf171508034449025ff985881/e2e/743a8dc0-4c9d-4dbf-9795-7434f645abdb.md"
$mdName1 = "3b55afe9-fca1-4f9d-b841-af26fca2fc20.md"
$mdName2 = "743a8dc0-4c9d-4dbf-9795-7434f645abdb.md"

# ---------------------------------------------------------------------------
# Overview sheet: refresh the rolled-up status column for both languages
# ---------------------------------------------------------------------------
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# Columns E/F on the Overview sheet grow to fit the new, longer status text.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn handback table
# ---------------------------------------------------------------------------
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("J2").Value = "3b55afe9-fca1-4f9d-b841-af26fca2fc20.783808eb53d60caf257a1f5621c198600b8a227c.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-06 09:01:52"

$zhcn.Range("J3").Value = "743a8dc0-4c9d-4dbf-9795-7434f645abdb.d548fde11194a939a491cb9b1bc79ede6d3f6d11.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-06 09:01:52"

# "Latest Target File" (column I) now links back to the source markdown file, same as column A.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $mdName1)
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $mdName1)
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $mdUrl2, [Type]::Missing, [Type]::Missing, $mdName2)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrl2, [Type]::Missing, [Type]::Missing, $mdName2)

# Column C (Status) grows to fit the new, longer status text; I/J grow for the file names.
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de handback table
# ---------------------------------------------------------------------------
$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("J2").Value = "3b55afe9-fca1-4f9d-b841-af26fca2fc20.783808eb53d60caf257a1f5621c198600b8a227c.de-de.xlf"
$dede.Range("K2").Value = "2016-09-06 09:02:27"

$dede.Range("J3").Value = "743a8dc0-4c9d-4dbf-9795-7434f645abdb.d548fde11194a939a491cb9b1bc79ede6d3f6d11.de-de.xlf"
$dede.Range("K3").Value = "2016-09-06 09:02:27"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $mdName1)
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $mdName1)
$dede.Hyperlinks.Add($dede.Range("A3"), $mdUrl2, [Type]::Missing, [Type]::Missing, $mdName2)
$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrl2, [Type]::Missing, [Type]::Missing, $mdName2)

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
